$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 7;  B = "/Envelope/Body/GetInfoByCityResponse/GetInfoByCityResult/NewDataSet/Table[2]/CITY";      C = "[A-Z a-z].*" },
    @{ Row = 8;  B = "/Envelope/Body/GetInfoByCityResponse/GetInfoByCityResult/NewDataSet/Table[2]/STATE";     C = "[A-Z]{2}" },
    @{ Row = 9;  B = "/Envelope/Body/GetInfoByCityResponse/GetInfoByCityResult/NewDataSet/Table[2]/ZIP";       C = "[0-9]{5}" },
    @{ Row = 10; B = "/Envelope/Body/GetInfoByCityResponse/GetInfoByCityResult/NewDataSet/Table[2]/AREA_CODE"; C = "[0-9]{3}" },
    @{ Row = 11; B = "/Envelope/Body/GetInfoByCityResponse/GetInfoByCityResult/NewDataSet/Table[2]/TIME_ZONE"; C = "[A-Z]{1}" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}

$ws.Range("B12").Select()
